$d = $word.ActiveDocument

$d.Content.Find.Execute("multi faceted", $true, $false, $false, $false, $false,
                         $true, 1, $false, "multi-faceted", 2)

$d.Content.Find.Execute("locical", $true, $false, $false, $false, $false,
                         $true, 1, $false, "logical", 2)
